$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Years_interpolation1")
$ws2 = $wb.Worksheets.Item("Years_interpolation2")

$years = @("year2007","year2010","year2015","year2020","year2025","year2030","year2035","year2040","year2045","year2050")

for ($i = 0; $i -lt $years.Length; $i++) {
    $col = 2 + $i
    $ws1.Cells.Item(2, $col).Value = $years[$i]
    $ws2.Cells.Item(2, $col).Value = $years[$i]
}

$ws2.Select()
$ws2.Range("B2:K2").Select()

$ws1.Select()
$ws1.Range("M8").Select()

Write-Output "Done"
